$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(22, 1).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(22, 1).Value = (Get-Date -Year 2020 -Month 3 -Day 6 -Hour 0 -Minute 0 -Second 0)

$ws.Cells.Item(23, 1).Value = 1111.223

$ws.Columns.Item(1).AutoFit() | Out-Null

$ws.Range("A24").Select()
